$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the precise timestamp value on row 11 (tiny precision correction)
$ws.Range("A11").Value = 45811.39382193287

# Append new row 12 with the latest price data
$ws.Range("A12").Value = 45812.39352998948
$ws.Range("B12").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C12").Value = "1Kg"
$ws.Range("D12").Value = "12,88€"

# Match the date/time number format and style used by the other rows in column A
$ws.Range("A12").NumberFormat = $ws.Range("A11").NumberFormat
